# Auto-generated edit script: update '想去人数' (F column) values across all 4 sheets
# per the diff (data refresh snapshot at commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 10529
$ws.Cells.Item(3, 6).Value = 246
$ws.Cells.Item(4, 6).Value = 1164
$ws.Cells.Item(5, 6).Value = 1034
$ws.Cells.Item(6, 6).Value = 821
$ws.Cells.Item(7, 6).Value = 265
$ws.Cells.Item(9, 6).Value = 338
$ws.Cells.Item(10, 6).Value = 1101
$ws.Cells.Item(13, 6).Value = 849
$ws.Cells.Item(14, 6).Value = 368
$ws.Cells.Item(15, 6).Value = 1795
$ws.Cells.Item(16, 6).Value = 18
$ws.Cells.Item(17, 6).Value = 898
$ws.Cells.Item(18, 6).Value = 803
$ws.Cells.Item(19, 6).Value = 532
$ws.Cells.Item(20, 6).Value = 767
$ws.Cells.Item(21, 6).Value = 868
$ws.Cells.Item(23, 6).Value = 267
$ws.Cells.Item(24, 6).Value = 85
$ws.Cells.Item(25, 6).Value = 594
$ws.Cells.Item(26, 6).Value = 608
$ws.Cells.Item(27, 6).Value = 104
$ws.Cells.Item(28, 6).Value = 317
$ws.Cells.Item(29, 6).Value = 993
$ws.Cells.Item(30, 6).Value = 37
$ws.Cells.Item(31, 6).Value = 482
$ws.Cells.Item(32, 6).Value = 146
$ws.Cells.Item(34, 6).Value = 212
$ws.Cells.Item(35, 6).Value = 537
$ws.Cells.Item(36, 6).Value = 1631
$ws.Cells.Item(37, 6).Value = 360
$ws.Cells.Item(39, 6).Value = 1388
$ws.Cells.Item(40, 6).Value = 401
$ws.Cells.Item(42, 6).Value = 44
$ws.Cells.Item(43, 6).Value = 76
$ws.Cells.Item(45, 6).Value = 66
$ws.Cells.Item(46, 6).Value = 66
$ws.Cells.Item(47, 6).Value = 34

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(5, 6).Value = 180
$ws.Cells.Item(11, 6).Value = 169

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 2144
$ws.Cells.Item(3, 6).Value = 601
$ws.Cells.Item(4, 6).Value = 520

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 2144
$ws.Cells.Item(3, 6).Value = 601
$ws.Cells.Item(5, 6).Value = 10529
$ws.Cells.Item(6, 6).Value = 246
$ws.Cells.Item(7, 6).Value = 1164
$ws.Cells.Item(8, 6).Value = 520
$ws.Cells.Item(9, 6).Value = 1034
$ws.Cells.Item(10, 6).Value = 821
$ws.Cells.Item(11, 6).Value = 180
$ws.Cells.Item(12, 6).Value = 338
$ws.Cells.Item(13, 6).Value = 1101
$ws.Cells.Item(16, 6).Value = 849
$ws.Cells.Item(17, 6).Value = 368
$ws.Cells.Item(18, 6).Value = 1795
$ws.Cells.Item(19, 6).Value = 18
$ws.Cells.Item(20, 6).Value = 898
$ws.Cells.Item(21, 6).Value = 803
$ws.Cells.Item(22, 6).Value = 532
$ws.Cells.Item(23, 6).Value = 767
$ws.Cells.Item(24, 6).Value = 868
$ws.Cells.Item(26, 6).Value = 267
$ws.Cells.Item(27, 6).Value = 85
$ws.Cells.Item(28, 6).Value = 594
$ws.Cells.Item(31, 6).Value = 608
$ws.Cells.Item(32, 6).Value = 104
$ws.Cells.Item(33, 6).Value = 317
$ws.Cells.Item(34, 6).Value = 993
$ws.Cells.Item(36, 6).Value = 37
$ws.Cells.Item(37, 6).Value = 482
$ws.Cells.Item(38, 6).Value = 146
$ws.Cells.Item(40, 6).Value = 360
$ws.Cells.Item(41, 6).Value = 1388
$ws.Cells.Item(42, 6).Value = 401
$ws.Cells.Item(45, 6).Value = 44
$ws.Cells.Item(46, 6).Value = 76
$ws.Cells.Item(47, 6).Value = 66
